# Insert a new weekly record for Apio (Terminal Hortofrutícola Agro Chillán)
# at row 52, shifting all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(52).Insert()

$ws.Cells.Item(52, 1).Value = 7
$ws.Cells.Item(52, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(52, 3).Value = "Ñuble"
$ws.Cells.Item(52, 4).Value = 44536
$ws.Cells.Item(52, 5).Value = 16
$ws.Cells.Item(52, 6).Value = 100112017
$ws.Cells.Item(52, 7).Value = "Apio"
$ws.Cells.Item(52, 8).Value = "Americana (o)"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 60
$ws.Cells.Item(52, 11).Value = 8000
$ws.Cells.Item(52, 12).Value = 8500
$ws.Cells.Item(52, 13).Value = 8250
$ws.Cells.Item(52, 14).Value = "`$/docena de matas"
$ws.Cells.Item(52, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(52, 16).Value = 1375
$ws.Cells.Item(52, 17).Value = 6
$ws.Cells.Item(52, 18).Value = "Hortaliza"
